$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<Name>_old" -> "<Name>_FV2310" (A1:J1)
#    and "<Name>_new" -> "<Name>_FV2404" (L1:U1). K1 ("diff") is unchanged.
# ---------------------------------------------------------------------------
$fv2310Headers = @("Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310")
$fv2404Headers = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U57 into an Excel Table ("Table1") with an AutoFilter, while
#    preserving the header row's existing formatting exactly (no new dxf).
#    We do this by stashing a copy of the header formatting in a scratch
#    row, clearing the header's direct formatting (so table creation does
#    not capture/clone it into a dxf), creating the table, then pasting the
#    original formatting back onto the header and wiping the scratch row.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A59:U59")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false

$scratchRange.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1) and select the bottom pane.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
